$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as literal text (matches existing rows, which are
# shared strings like "2020-06-01" rather than real dates). Force text entry
# by setting the number format to Text before assigning the value, otherwise
# the "YYYY-MM-DD" string gets auto-converted into a date serial number.
# Reset the style back to Normal afterwards so no extra cell formatting is
# left behind on the new cells.
$ws.Range("A62").NumberFormat = "@"
$ws.Range("A62").Value = "2020-07-31"
$ws.Range("A62").Style = "Normal"
$ws.Range("B62").Value = 424637
$ws.Range("C62").Value = 469629
$ws.Range("D62").Value = 90022
$ws.Range("E62").Value = 46688
$ws.Range("F62").Value = 27.16

$ws.Range("A63").NumberFormat = "@"
$ws.Range("A63").Value = "2020-08-01"
$ws.Range("A63").Style = "Normal"
$ws.Range("B63").Value = 434193
$ws.Range("C63").Value = 477733
$ws.Range("D63").Value = 87771
$ws.Range("E63").Value = 47472
$ws.Range("F63").Value = 27.02
